$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "98.068.70"
$ws.Range("E2").Value = "  -0.46%  "
$ws.Range("D3").Value = "3.395.92"
$ws.Range("E3").Value = "  +0.82%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "254.17"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -1.48%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "681.92"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +1.65%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.44"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -7.66%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.429"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -7.02%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1.05"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -4.31%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.999"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +0.01%  "
$ws.Range("D11").Value = "3.390.11"
$ws.Range("E11").Value = "  +0.71%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.215"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +2.23%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "41.45"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -2.44%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.25"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +11.01%  "
$ws.Range("D15").Value = "97.703.50"
$ws.Range("E15").Value = "  -0.55%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000263"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -1.95%  "
$ws.Range("D17").Value = "4.032.00"
$ws.Range("E17").Value = "  +0.70%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "8.88"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +16.02%  "
$ws.Range("D19").Value = "3.431.98"
$ws.Range("E19").Value = "  +1.79%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.566"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +25.36%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.38"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +2.70%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.96"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +3.38%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.41"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -4.58%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "505.31"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -4.81%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0000203"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -4.75%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.54"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +4.18%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "99.81"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -2.83%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "12.63"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -0.53%  "
$ws.Range("D29").Value = "3.593.23"
$ws.Range("E29").Value = "  +1.13%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.149"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -0.58%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "11.47"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +2.39%  "
$ws.Range("E32").Value = "  +0.01%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.194"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +2.50%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.62"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +22.83%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.00"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +0.06%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.570"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +5.14%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "29.55"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -0.94%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "7.92"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +0.18%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.51"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +12.32%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "524.21"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -0.62%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.152"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -5.77%  "
$ws.Range("E42").Value = "  -0.09%  "
$ws.Range("B43").Value = "WhiteBITCoin"
$ws.Range("C43").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "24.71"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +0.01%  "
$ws.Range("B44").Value = "ARBITRUM"
$ws.Range("C44").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.864"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +2.55%  "
$ws.Range("B45").Value = "MantraDAO"
$ws.Range("C45").Value = "https://coinranking.com/coin/cTdD8lD-6+mantradao-om"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.79"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -0.19%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0433"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -4.74%  "
$ws.Range("B47").Value = "Cosmos"
$ws.Range("C47").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.91"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +11.68%  "
$ws.Range("B48").Value = "ImmutableX"
$ws.Range("C48").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.72"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +12.09%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "5.72"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +10.27%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "55.86"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +10.04%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "3.19"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -5.56%  "
